# Applies updated dSF (column F) values as reflected in the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    7  = -5
    8  = 1
    9  = 3
    12 = -1
    13 = 2
    14 = -4
    23 = -3
    25 = -3
    26 = -2
    27 = -2
    29 = -1
    31 = 5
    32 = 4
    33 = 4
    34 = -1
    36 = 0
    43 = 3
    47 = 1
    48 = -5
    49 = -8
    50 = 6
    51 = -5
    54 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
